$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-Paragraph($index, $xml) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $null = $r.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Paragraph 50: ORTHOGRAPHIC_CAMERA : CAMERA
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:pPr><w:ind w:left=`"720`"/></w:pPr>" +
         "<w:r><w:t>ORTHOGRAPHIC_</w:t></w:r>" +
         "<w:proofErr w:type=`"gramStart`"/>" +
         "<w:r><w:t>CAMERA</w:t></w:r>" +
         "<w:r><w:t xml:space=`"preserve`"> :</w:t></w:r>" +
         "<w:proofErr w:type=`"gramEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> CAMERA</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 50 $xml

# ---------------------------------------------------------------------------
# Paragraph 49: PERSPECTIVE_CAMERA : CAMERA  -> + new FieldOfView paragraph
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:pPr><w:ind w:left=`"720`"/></w:pPr>" +
         "<w:r><w:t>PERSPECTIVE_</w:t></w:r>" +
         "<w:proofErr w:type=`"gramStart`"/>" +
         "<w:r><w:t>CAMERA</w:t></w:r>" +
         "<w:r><w:t xml:space=`"preserve`"> :</w:t></w:r>" +
         "<w:proofErr w:type=`"gramEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> CAMERA</w:t></w:r>" +
       "</w:p>" +
       "<w:p $wNs>" +
         "<w:pPr>" +
           "<w:pStyle w:val=`"ListParagraph`"/>" +
           "<w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr>" +
         "</w:pPr>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:r><w:t>FieldOfView</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> (float) – [Radians]</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 49 $xml

# ---------------------------------------------------------------------------
# Paragraph 48: FarPlane (float)
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:pPr>" +
           "<w:pStyle w:val=`"ListParagraph`"/>" +
           "<w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr>" +
         "</w:pPr>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:r><w:t>FarPlane</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> (float)</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 48 $xml

# ---------------------------------------------------------------------------
# Paragraph 47: UpDirection (Vector3)
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:pPr>" +
           "<w:pStyle w:val=`"ListParagraph`"/>" +
           "<w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr>" +
         "</w:pPr>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:r><w:t>UpDirection</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> (Vector3)</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 47 $xml

# ---------------------------------------------------------------------------
# Paragraph 44: FocalLength (float)
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:pPr>" +
           "<w:pStyle w:val=`"ListParagraph`"/>" +
           "<w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr>" +
         "</w:pPr>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:r><w:t>FocalLength</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> (float)</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 44 $xml

# ---------------------------------------------------------------------------
# Paragraph 36: VertexD (Vector3)
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:pPr>" +
           "<w:pStyle w:val=`"ListParagraph`"/>" +
           "<w:ind w:firstLine=`"720`"/>" +
         "</w:pPr>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:r><w:t>VertexD</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> (Vector3)</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 36 $xml

# ---------------------------------------------------------------------------
# Paragraph 35: RECTANGLE : TRIANGLE
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:pPr><w:pStyle w:val=`"ListParagraph`"/></w:pPr>" +
         "<w:proofErr w:type=`"gramStart`"/>" +
         "<w:r><w:t>RECTANGLE</w:t></w:r>" +
         "<w:r><w:t xml:space=`"preserve`"> :</w:t></w:r>" +
         "<w:proofErr w:type=`"gramEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> TRIANGLE</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 35 $xml

# ---------------------------------------------------------------------------
# Paragraph 33: VertexC (Vector3) + trailing tab
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:pPr><w:pStyle w:val=`"ListParagraph`"/></w:pPr>" +
         "<w:r><w:tab/></w:r>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:r><w:t>VertexC</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> (Vector3)</w:t></w:r>" +
         "<w:r><w:tab/></w:r>" +
       "</w:p>"
Replace-Paragraph 33 $xml

# ---------------------------------------------------------------------------
# Paragraph 32: VertexB (Vector3),
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:pPr><w:pStyle w:val=`"ListParagraph`"/></w:pPr>" +
         "<w:r><w:tab/></w:r>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:r><w:t>VertexB</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> (Vector3),</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 32 $xml

# ---------------------------------------------------------------------------
# Paragraph 31: VertexA (Vector3),
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:pPr><w:pStyle w:val=`"ListParagraph`"/></w:pPr>" +
         "<w:r><w:tab/></w:r>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:r><w:t>VertexA</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> (Vector3),</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 31 $xml

# ---------------------------------------------------------------------------
# Paragraph 26: Height  and Width (floats)
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:pPr>" +
           "<w:pStyle w:val=`"ListParagraph`"/>" +
           "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr>" +
           "<w:rPr><w:b/><w:bCs/></w:rPr>" +
         "</w:pPr>" +
         "<w:proofErr w:type=`"gramStart`"/>" +
         "<w:r><w:t xml:space=`"preserve`">Height </w:t></w:r>" +
         "<w:r><w:t xml:space=`"preserve`"> and</w:t></w:r>" +
         "<w:proofErr w:type=`"gramEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> Width </w:t></w:r>" +
         "<w:r><w:t>(floats)</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 26 $xml

# ---------------------------------------------------------------------------
# Paragraph 25: SCENE (namespace: SceneObjects)
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:r><w:t>SCENE</w:t></w:r>" +
         "<w:r><w:t xml:space=`"preserve`"> (namespace: </w:t></w:r>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:r><w:t>SceneObjects</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
         "<w:r><w:t>)</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 25 $xml

# ---------------------------------------------------------------------------
# Paragraph 23: X,y,z  -> + two new italic paragraphs (Float Magnitude / Vector3 Normalize)
# ---------------------------------------------------------------------------
$xml = "<w:p $wNs>" +
         "<w:r><w:tab/></w:r>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:proofErr w:type=`"gramStart`"/>" +
         "<w:r><w:t>X,y</w:t></w:r>" +
         "<w:proofErr w:type=`"gramEnd`"/>" +
         "<w:r><w:t>,z</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
       "</w:p>" +
       "<w:p $wNs>" +
         "<w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr>" +
         "<w:r><w:tab/></w:r>" +
         "<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=`"preserve`">Float </w:t></w:r>" +
         "<w:proofErr w:type=`"gramStart`"/>" +
         "<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Magnitude(</w:t></w:r>" +
         "<w:proofErr w:type=`"gramEnd`"/>" +
         "<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>)</w:t></w:r>" +
       "</w:p>" +
       "<w:p $wNs>" +
         "<w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr>" +
         "<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:tab/><w:t xml:space=`"preserve`">Vector3 </w:t></w:r>" +
         "<w:proofErr w:type=`"gramStart`"/>" +
         "<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Normalize(</w:t></w:r>" +
         "<w:proofErr w:type=`"gramEnd`"/>" +
         "<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>)</w:t></w:r>" +
       "</w:p>"
Replace-Paragraph 23 $xml
